$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the summary rows (tunnit yht. / target / suoritettu(%)) down by 4 rows
# (from 108:110 to 112:114), opening up rows 104:109 for new log entries and
# leaving rows 110:111 blank, matching the target layout.
$ws.Rows("108:111").Insert()

# --- New time-tracking entries (rows 104-109) ---

# Row 104: new date 27.1.2022 (44588). Clone the date formatting (s="5",
# m/d/yyyy, centered) from an existing date cell via PasteSpecial-Formats
# (xlPasteFormats = -4122) so the new cell reuses the same cell style.
$ws.Range("A97").Copy()
$ws.Range("A104").PasteSpecial(-4122)
$ws.Range("A104").Value = 44588
$ws.Range("B104").Value = 1
$ws.Range("C104").Value = "tustustuminen sass react dokumentointiin"
$ws.Range("D104").Value = "client"

$ws.Range("B106").Value = 2
$ws.Range("C106").Value = "form submit ja reset komponettien erottelu ja scss, SignUpFormCntr jotta keskustelu backending kanssa erillään, extractErrorMsg"
$ws.Range("D106").Value = "client"

$ws.Range("B105").Value = 1
$ws.Range("C105").Value = ".scss testaus, import, extend jne, Constants.scss aloitettu"
$ws.Range("D105").Value = "client"

$ws.Range("B107").Value = 2
$ws.Range("C107").Value = "loginForm refaktoroitu, scss ja LoginFormCntr"
$ws.Range("D107").Value = "client"

$ws.Range("B108").Value = 3
$ws.Range("C108").Value = "analyzerForm, refaktorointia ja scss, visuaalista tuunausta, "
$ws.Range("D108").Value = "client"

$ws.Range("B109").Value = 1
$ws.Range("C109").Value = "analyzerInfo scss ja pientä tuunausta"
$ws.Range("D109").Value = "client"

# Row 106: new date 28.1.2022 (44589), same style cloning as above
$ws.Range("A97").Copy()
$ws.Range("A106").PasteSpecial(-4122)
$ws.Range("A106").Value = 44589

# --- Update summary formulas for the new data range ---
$ws.Range("B112").Formula = "=SUM(B2:B109)"
$ws.Range("B114").Formula = "=B112/B113*100"

# --- Restore view state (scroll position / selection / window placement) ---
$excel.ActiveWindow.ScrollRow = 94
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C109").Select()

$wb.Windows.Item(1).Left = 3720
$wb.Windows.Item(1).Top = 1515
